$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "25.913.71"
Set-TextValue "E2" "  -0.34%  "
Set-TextValue "D3" "1.623.08"
Set-TextValue "E4" "  +0.20%  "
Set-TextValue "D5" "213.50"
Set-TextValue "E5" "  -1.09%  "
Set-TextValue "E6" "  -0.90%  "
Set-TextValue "E8" "  -2.19%  "
Set-TextValue "E9" "  -3.68%  "
Set-TextValue "D10" "18.26"
Set-TextValue "E10" "  -6.56%  "
Set-TextValue "D11" "0.0787"
Set-TextValue "E11" "  -1.19%  "
Set-TextValue "D12" "1.849.02"
Set-TextValue "E12" "  -1.09%  "
Set-TextValue "D13" "1.623.12"
Set-TextValue "E13" "  -1.37%  "
Set-TextValue "D15" "0.523"
Set-TextValue "E15" "  -3.83%  "
Set-TextValue "D16" "25.909.51"
Set-TextValue "E16" "  -0.75%  "
Set-TextValue "D17" "61.09"
Set-TextValue "E17" "  -3.57%  "
Set-TextValue "D18" "0.0₃0732"
Set-TextValue "E18" "  -3.99%  "
Set-TextValue "E19" "  +0.21%  "
Set-TextValue "D20" "191.76"
Set-TextValue "E20" "  -1.18%  "
Set-TextValue "E21" "  -3.09%  "
Set-TextValue "D22" "9.56"
Set-TextValue "E22" "  -3.63%  "
Set-TextValue "E23" "  -2.22%  "
Set-TextValue "D24" "0.132"
Set-TextValue "E24" "  +0.83%  "
Set-TextValue "D25" "143.81"
Set-TextValue "E25" "  +0.60%  "
Set-TextValue "E26" "  +0.28%  "
Set-TextValue "E27" "  -3.04%  "
Set-TextValue "D28" "6.72"
Set-TextValue "E28" "  -2.21%  "
Set-TextValue "D29" "15.13"
Set-TextValue "E29" "  -2.49%  "
Set-TextValue "E30" "  -1.38%  "
Set-TextValue "D31" "0.0482"
Set-TextValue "E31" "  -2.65%  "
Set-TextValue "E32" "  -4.26%  "
Set-TextValue "E33" "  -5.53%  "
Set-TextValue "E34" "  -2.96%  "
Set-TextValue "E35" "  -2.46%  "
Set-TextValue "D36" "1.116.79"
Set-TextValue "E36" "  -1.20%  "
Set-TextValue "E37" "  -6.53%  "
Set-TextValue "E38" "  -1.32%  "
Set-TextValue "D39" "0.516"
Set-TextValue "E39" "  -4.41%  "
Set-TextValue "E40" "  -2.40%  "
Set-TextValue "D41" "97.88"
Set-TextValue "E41" "  -1.18%  "
Set-TextValue "D42" "0.766"
Set-TextValue "E42" "  -3.91%  "
Set-TextValue "D43" "1.759.06"
Set-TextValue "E43" "  -1.10%  "
Set-TextValue "E44" "  -5.80%  "
Set-TextValue "E45" "  -1.64%  "
Set-TextValue "E46" "  +1.66%  "
Set-TextValue "D47" "54.31"
Set-TextValue "E47" "  -3.89%  "
Set-TextValue "E48" "  -1.44%  "
Set-TextValue "D49" "0.412"
Set-TextValue "E49" "  -0.44%  "
Set-TextValue "E50" "  +0.26%  "
Set-TextValue "D51" "7.48"
Set-TextValue "E51" "  -3.23%  "
